$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 105 ---
$ws.Range("A105:B105").NumberFormat = "@"

$ws.Range("A105").Value = "14369246"
$ws.Range("B105").Value = "2025-08-07"
$ws.Range("C105").Value = "Anastasia Potapova"
$ws.Range("D105").Value = "Laura Siegemund"
$ws.Range("E105").Value = "Gana Anastasia Potapova"
$ws.Range("F105").Value = 1.8

$ws.Range("A105:B105").Style = "Normal"

# --- Row 106 ---
$ws.Range("A106:B106").NumberFormat = "@"

$ws.Range("A106").Value = "14311060"
$ws.Range("B106").Value = "2025-08-08"
$ws.Range("C106").Value = "Santiago Rodriguez Taverna"
$ws.Range("D106").Value = "Dusan Lajovic"
$ws.Range("E106").Value = "Gana Santiago Rodriguez Taverna"
$ws.Range("F106").Value = 4

$ws.Range("A106:B106").Style = "Normal"
